# Daily update at 8 AM UTC
# Appends the next day's win-count row to the "Wins Over Time" tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the current last used row (before this update it is row 68).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# New day's data: date serial, Chase, Bryce, Zach win counts.
$ws.Cells.Item($newRow, 1).Value = 46018
$ws.Cells.Item($newRow, 2).Value = 153
$ws.Cells.Item($newRow, 3).Value = 162
$ws.Cells.Item($newRow, 4).Value = 151

# Keep the date column formatted/styled like the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
